$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update the "Date" property value (row 8, column B) ---
$ws.Cells.Item(8, 2).Value = "2024-10-02T15:04:17+00:00"

# --- Update the "Contact" property value (row 10, column B) ---
$ws.Cells.Item(10, 2).Value = "Ferlab.bio (http://example.org/example-publisher)"

# --- Insert a new "Jurisdiction" property row right after "Contact" (row 11) ---
# This pushes the existing rows 11-21 (Description ... Count) down to 12-22.
$ws.Rows.Item(11).Insert()

# Match the formatting of the other data rows (border + top-aligned wrapped text)
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
